$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -13.2051
$ws.Range("B12").Value = 5.001199999999997
$ws.Range("C23").Value = -11.9973
$ws.Range("B27").Value = 6.212900000000003
$ws.Range("C28").Value = -13.3051
$ws.Range("B32").Value = 6.605199999999996
$ws.Range("C32").Value = -11.7413
$ws.Range("C34").Value = -11.74140000000001
$ws.Range("B36").Value = 9.197100000000008
$ws.Range("B38").Value = 5.316999999999997
$ws.Range("C42").Value = -12.23730000000001
$ws.Range("B46").Value = 5.949100000000004
$ws.Range("C49").Value = -13.7572
$ws.Range("B54").Value = 5.932300000000002
$ws.Range("C54").Value = -12.661
$ws.Range("B55").Value = 5.250599999999996
$ws.Range("B56").Value = 4.599999999999994
$ws.Range("B67").Value = 5.207599999999997
$ws.Range("B69").Value = 5.269499999999997
$ws.Range("B72").Value = 5.079500000000003
$ws.Range("C78").Value = -12.9865
$ws.Range("C80").Value = -13.5343
$ws.Range("B83").Value = 5.922899999999998
$ws.Range("B86").Value = 5.154900000000001
$ws.Range("B91").Value = 5.169499999999995
$ws.Range("B93").Value = 5.1879
$ws.Range("C97").Value = -11.7142
$ws.Range("B99").Value = 6.252799999999999
$ws.Range("C99").Value = -11.87320000000001
$ws.Range("C101").Value = -12.7525
$ws.Range("B104").Value = 9.4803
